# Update Mappings 22 Ontologies
# Adds a new "OSMO_DEF" column (F) to Sheet1, header styled like the other
# header cells, and filled with the literal string "[]" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the style used by the other header cells
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new column for each data row (rows 2-9) with the literal "[]"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
